# Update rows 2-9 of Sheet1 with the new TPM-derived values (per commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$row2 = @("FAPs", "Wnt5a", "Fzd3", "ECs", 3, 1, 8.775006, 26.325018, 0.9920592728348052, 0.9920592728348053, 2, 0.6666666666666666, 0.3251496666666667, 0.975449, 0.07121046526627427, 0.07121046526627427, 2.853190275898, 25.678712483082, 0.0706450023902882, 0.07064500239028822)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "2").Value = $row2[$i] }

$row3 = @("FAPs", "Wnt5a", "Fzd3", "FAPs", 3, 1, 8.775006, 26.325018, 0.9920592728348052, 0.9920592728348053, 3, 1, 1.007819666666667, 3.023459, 0.2207208394324094, 0.2207208394324094, 8.843623621917999, 79.592612597262, 0.2189681554668039, 0.2189681554668039)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "3").Value = $row3[$i] }

$row4 = @("FAPs", "Wnt5a", "Fzd3", "MuSCs", 3, 1, 8.775006, 26.325018, 0.9920592728348052, 0.9920592728348053, 3, 1, 3.226895, 9.680685, 0.7067166842615477, 0.7067166842615475, 28.31602298637, 254.84420687733, 0.7011048398887356, 0.7011048398887356)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "4").Value = $row4[$i] }

$row5 = @("FAPs", "Wnt5a", "Fzd3", "Resolving-Mac", 3, 1, 8.775006, 26.325018, 0.9920592728348052, 0.9920592728348053, 1, 0.3333333333333333, 0.006173333333333333, 0.01852, 0.001352011039768762, 0.001352011039768762, 0.05417103703999999, 0.4875393333599999, 0.001341275088977627, 0.001341275088977627)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "5").Value = $row5[$i] }

$row6 = @("MuSCs", "Wnt5a", "Fzd3", "ECs", 1, 0.3333333333333333, 0.07023766666666667, 0.210713, 0.007940727165194733, 0.007940727165194734, 2, 0.6666666666666666, 0.3251496666666667, 0.975449, 0.07121046526627427, 0.07121046526627427, 0.02283775390411111, 0.205539785137, 0.00056546287598606, 0.0005654628759860601)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "6").Value = $row6[$i] }

$row7 = @("MuSCs", "Wnt5a", "Fzd3", "FAPs", 1, 0.3333333333333333, 0.07023766666666667, 0.210713, 0.007940727165194733, 0.007940727165194734, 3, 1, 1.007819666666667, 3.023459, 0.2207208394324094, 0.2207208394324094, 0.07078690180744444, 0.6370821162670001, 0.001752683965605518, 0.001752683965605518)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "7").Value = $row7[$i] }

$row8 = @("MuSCs", "Wnt5a", "Fzd3", "MuSCs", 1, 0.3333333333333333, 0.07023766666666667, 0.210713, 0.007940727165194733, 0.007940727165194734, 3, 1, 3.226895, 9.680685, 0.7067166842615477, 0.7067166842615475, 0.2266495753783334, 2.039846178405, 0.00561184437281202, 0.00561184437281202)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "8").Value = $row8[$i] }

$row9 = @("MuSCs", "Wnt5a", "Fzd3", "Resolving-Mac", 1, 0.3333333333333333, 0.07023766666666667, 0.210713, 0.007940727165194733, 0.007940727165194734, 1, 0.3333333333333333, 0.006173333333333333, 0.01852, 0.001352011039768762, 0.001352011039768762, 0.0004336005288888889, 0.00390240476, 0.00001073595079113499, 0.00001073595079113499)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "9").Value = $row9[$i] }

# Row 10 (ECs/MuSCs target) no longer exists in the updated dataset; remove it and shrink the used range to A1:T9
$ws.Rows.Item(10).Delete()
